# Card drawing bug fixes:
#  - "Cloned" cards (sub_type action cards whose rows were duplicated)
#    had their F-column (pop_change_rate) values scaled up so the
#    background effect reflects the full stacked amount instead of the
#    single-card amount that was only showing the description.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Negative "population cost" cards: -500 -> -1300 for rows 8-16
$ws.Range("F8:F16").Value = -1300

# Negative "population cost" cards: -500 -> -900 for rows 17-23
$ws.Range("F17:F23").Value = -900

# Positive "population gain" cards: scaled up by 10x for rows 24-28
$ws.Range("F24").Value = 6000
$ws.Range("F25").Value = 5000
$ws.Range("F26").Value = 7000
$ws.Range("F27").Value = 3000
$ws.Range("F28").Value = 6000

# Restore the author's last on-screen selection/zoom state from the
# editing session.
$ws.Range("F17:F23").Select()
$excel.ActiveWindow.Zoom = 100
